$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for employee code column ("EmpCd" -> "Employee Code") and widen the column
$ws.Range("A1").Value = "Employee Code"
$ws.Columns.Item(1).ColumnWidth = 14.21875

# The stray "Payable"/"NHrs" headers are replaced with what used to be the W.OT/H.OT headers
$ws.Range("E1").Value = "W.OT"
$ws.Range("F1").Value = "H.OT"

# Remove the now-unused trailing columns from the header row
$ws.Range("G1").Clear()
$ws.Range("H1").Clear()
$ws.Range("I1").ClearContents()

# Row 2 sample data no longer needs the trailing three columns either
$ws.Range("G2:I2").ClearContents()

# Update selection to match what was saved
$ws.Range("C7").Select()

# Match the saved window geometry as closely as the object model allows
$win = $excel.ActiveWindow
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12456
